$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.155.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "'3.973.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'542.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.34%  "
$ws.Range("D6").Value = "'149.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").Value = "'3.966.94"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("D8").Value = "'0.689"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.10%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.744"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("D11").Value = "'0.167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.39%  "
$ws.Range("D12").Value = "'56.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +18.90%  "
$ws.Range("D13").Value = "'0.0000318"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "'10.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("D15").Value = "'4.614.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "'3.981.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'20.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'13.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").Value = "'71.254.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "'428.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").Value = "'97.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.79%  "
$ws.Range("D24").Value = "'3.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "'4.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.18%  "
$ws.Range("D26").Value = "'14.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").Value = "'11.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "'10.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").Value = "'3.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +17.66%  "
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "'36.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").Value = "'7.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.90%  "
$ws.Range("D33").Value = "'50.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +20.51%  "
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "'13.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "'681.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "'65.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").Value = "'0.438"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "'0.0₃0820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.150"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'3.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'0.0484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  -6.98%  "
$ws.Range("D47").Value = "'2.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'9.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.16%  "
$ws.Range("D49").Value = "'3.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.67%  "
$ws.Range("D50").Value = "'3.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").Value = "'0.000272"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.91%  "
